$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.482.52"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "'1.839.98"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +1.18%  "
$ws.Range("D5").Value = "'314.89"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("D7").Value = "'0.4755"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("D8").Value = "'0.3702"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("D9").Value = "'0.07473"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("D10").Value = "'0.8863"
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("D11").Value = "'20.52"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "'1.884.51"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").Value = "'0.07356"
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("D14").Value = "'5.458"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "'93.29"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "'6.592"
$ws.Range("D17").Value = "'1.013"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "'0.000008830"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").Value = "'1.012"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "'14.82"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").Value = "'27.497.85"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").Value = "'5.330"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "'10.71"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").Value = "'2.103.23"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("D25").Value = "'1.894"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'152.27"
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("D27").Value = "'18.63"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("D28").Value = "'2.159"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'5.250"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "'118.04"
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("D31").Value = "'0.09007"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").Value = "'0.7576"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "'1.182"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").Value = "'4.560"
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("D35").Value = "'2.966"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").Value = "'0.05345"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").Value = "'0.01959"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").Value = "'3.004"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "'7.335"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("E42").Value = "  +4.65%  "
$ws.Range("D43").Value = "'0.5347"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("D44").Value = "'0.1663"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").Value = "'0.4929"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "'10.53"
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("D48").Value = "'1.013"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").Value = "'104.92"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").Value = "'1.682"
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "'0.06317"
$ws.Range("E51").Value = "  +0.44%  "